$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: give the new columns L and M the correct cell format (style) on every
# row that holds data in D:K, by copying the format already present in column K.
$dataRows = @(7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102)
foreach ($r in $dataRows) {
    $ws.Range("K$r").Copy() | Out-Null
    $ws.Range("L$r`:M$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Step 2: write the final values (two new quarters inserted, some quarters restated)
# for every row, columns D through M.
$row = New-Object 'object[,]' 1,10
$row[0,0] = 43465
$row[0,1] = 43373
$row[0,2] = 43281
$row[0,3] = 43190
$row[0,4] = 43100
$row[0,5] = 43008
$row[0,6] = 42916
$row[0,7] = 42825
$row[0,8] = 42735
$row[0,9] = 42643
$ws.Range("D7:M7").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 431500
$row[0,1] = 412300
$row[0,2] = 416200
$row[0,3] = 421000
$row[0,4] = 400300
$row[0,5] = 405900
$row[0,6] = 410100
$row[0,7] = 425700
$row[0,8] = 394600
$row[0,9] = 347100
$ws.Range("D8:M8").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 206400
$row[0,1] = 202700
$row[0,2] = 197100
$row[0,3] = 188000
$row[0,4] = 187200
$row[0,5] = 188700
$row[0,6] = 191400
$row[0,7] = 192100
$row[0,8] = 182100
$row[0,9] = 153400
$ws.Range("D9:M9").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 225100
$row[0,1] = 209600
$row[0,2] = 219200
$row[0,3] = 233000
$row[0,4] = 213100
$row[0,5] = 217200
$row[0,6] = 218700
$row[0,7] = 233600
$row[0,8] = 212500
$row[0,9] = 193700
$ws.Range("D10:M10").Value = $row
$ws.Range("D11:M11").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = "NA"
$row[0,1] = "NA"
$row[0,2] = "NA"
$row[0,3] = "NA"
$row[0,4] = "NA"
$row[0,5] = "NA"
$row[0,6] = "NA"
$row[0,7] = "NA"
$row[0,8] = "NA"
$row[0,9] = "NA"
$ws.Range("D12:M12").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D13:M13").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 13600
$row[0,1] = 6400
$row[0,2] = 10800
$row[0,3] = 3800
$row[0,4] = 900
$row[0,5] = 19600
$row[0,6] = 109000
$row[0,7] = 3100
$row[0,8] = 6100
$row[0,9] = 1400
$ws.Range("D14:M14").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 46900
$row[0,1] = 44200
$row[0,2] = 46000
$row[0,3] = 43200
$row[0,4] = 43500
$row[0,5] = 42700
$row[0,6] = 46800
$row[0,7] = 45300
$row[0,8] = 42600
$row[0,9] = 36200
$ws.Range("D15:M15").Value = $row
$ws.Range("D16:M16").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 359500
$row[0,1] = 357700
$row[0,2] = 278400
$row[0,3] = 313200
$row[0,4] = 200800
$row[0,5] = 349900
$row[0,6] = 441900
$row[0,7] = 335100
$row[0,8] = 322300
$row[0,9] = 274000
$ws.Range("D17:M17").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 72000
$row[0,1] = 54600
$row[0,2] = 137800
$row[0,3] = 107800
$row[0,4] = 199500
$row[0,5] = 56000
$row[0,6] = -31800
$row[0,7] = 90600
$row[0,8] = 72300
$row[0,9] = 73100
$ws.Range("D18:M18").Value = $row
$ws.Range("D19:M19").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = -61300
$row[0,1] = -28900
$row[0,2] = -23900
$row[0,3] = -14800
$row[0,4] = -20000
$row[0,5] = -31300
$row[0,6] = -30100
$row[0,7] = -34500
$row[0,8] = -35300
$row[0,9] = -34900
$ws.Range("D20:M20").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 57600
$row[0,1] = 69900
$row[0,2] = 159900
$row[0,3] = 136200
$row[0,4] = 223000
$row[0,5] = 67300
$row[0,6] = -15200
$row[0,7] = 101400
$row[0,8] = 79600
$row[0,9] = 74500
$ws.Range("D21:M21").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D22:M22").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 10700
$row[0,1] = 25700
$row[0,2] = 113900
$row[0,3] = 93000
$row[0,4] = 179500
$row[0,5] = 24700
$row[0,6] = -62000
$row[0,7] = 56100
$row[0,8] = 37000
$row[0,9] = 38200
$ws.Range("D23:M23").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = -2400
$row[0,1] = 600
$row[0,2] = 14800
$row[0,3] = 10900
$row[0,4] = 48200
$row[0,5] = 2400
$row[0,6] = -11800
$row[0,7] = 10700
$row[0,8] = -4100
$row[0,9] = 4800
$ws.Range("D24:M24").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D25:M25").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 13200
$row[0,1] = 25100
$row[0,2] = 99100
$row[0,3] = 82100
$row[0,4] = 131300
$row[0,5] = 22300
$row[0,6] = -50200
$row[0,7] = 45400
$row[0,8] = 41100
$row[0,9] = 33400
$ws.Range("D26:M26").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 8900
$row[0,1] = 14700
$row[0,2] = 82700
$row[0,3] = 51200
$row[0,4] = 114800
$row[0,5] = 11800
$row[0,6] = -25700
$row[0,7] = 19900
$row[0,8] = 20400
$row[0,9] = 8300
$ws.Range("D27:M27").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D28:M28").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = "NA"
$row[0,2] = "NA"
$row[0,3] = "NA"
$row[0,4] = -85300
$row[0,5] = "NA"
$row[0,6] = "NA"
$row[0,7] = "NA"
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D29:M29").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D30:M30").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D31:M31").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 61300
$row[0,1] = 28900
$row[0,2] = 23900
$row[0,3] = 14800
$row[0,4] = 20000
$row[0,5] = 31300
$row[0,6] = 30100
$row[0,7] = 34500
$row[0,8] = 35300
$row[0,9] = 34900
$ws.Range("D32:M32").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 8900
$row[0,1] = 14700
$row[0,2] = 82700
$row[0,3] = 51200
$row[0,4] = 29500
$row[0,5] = 11800
$row[0,6] = -25700
$row[0,7] = 19900
$row[0,8] = 20400
$row[0,9] = 8300
$ws.Range("D33:M33").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D34:M34").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 8900
$row[0,1] = 14700
$row[0,2] = 82700
$row[0,3] = 51200
$row[0,4] = 29500
$row[0,5] = 11800
$row[0,6] = -25700
$row[0,7] = 19900
$row[0,8] = 20400
$row[0,9] = 8300
$ws.Range("D35:M35").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 43465
$row[0,1] = 43373
$row[0,2] = 43281
$row[0,3] = 43190
$row[0,4] = 43100
$row[0,5] = 43008
$row[0,6] = 42916
$row[0,7] = 42825
$row[0,8] = 42735
$row[0,9] = 42643
$ws.Range("D38:M38").Value = $row
$ws.Range("D39:M39").ClearContents() | Out-Null
$ws.Range("D40:M40").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 114600
$row[0,1] = 110600
$row[0,2] = 108400
$row[0,3] = 179200
$row[0,4] = 231500
$row[0,5] = 222400
$row[0,6] = 125300
$row[0,7] = 119400
$row[0,8] = 133800
$row[0,9] = 96300
$ws.Range("D41:M41").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D42:M42").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 51400
$row[0,1] = 52400
$row[0,2] = 44100
$row[0,3] = 45500
$row[0,4] = 49000
$row[0,5] = 48800
$row[0,6] = 50700
$row[0,7] = 46300
$row[0,8] = 51200
$row[0,9] = 35100
$ws.Range("D43:M43").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 14900
$row[0,1] = 12700
$row[0,2] = 12700
$row[0,3] = 11600
$row[0,4] = 12600
$row[0,5] = 11500
$row[0,6] = 11100
$row[0,7] = 10900
$row[0,8] = 12000
$row[0,9] = 9400
$ws.Range("D44:M44").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 81100
$row[0,1] = 68200
$row[0,2] = 46900
$row[0,3] = 50100
$row[0,4] = 48500
$row[0,5] = 308100
$row[0,6] = 58600
$row[0,7] = 59100
$row[0,8] = 52900
$row[0,9] = 53100
$ws.Range("D45:M45").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 262000
$row[0,1] = 243800
$row[0,2] = 212200
$row[0,3] = 286400
$row[0,4] = 341600
$row[0,5] = 590800
$row[0,6] = 245700
$row[0,7] = 235700
$row[0,8] = 249800
$row[0,9] = 193900
$ws.Range("D46:M46").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 8900
$row[0,1] = 9200
$row[0,2] = 9500
$row[0,3] = 9700
$row[0,4] = 10100
$row[0,5] = 10200
$row[0,6] = 10000
$row[0,7] = 10200
$row[0,8] = 10600
$row[0,9] = 10700
$ws.Range("D47:M47").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 3206100
$row[0,1] = 3059700
$row[0,2] = 2939200
$row[0,3] = 2814500
$row[0,4] = 5261400
$row[0,5] = 2667400
$row[0,6] = 2644500
$row[0,7] = 2611900
$row[0,8] = 2601800
$row[0,9] = 2304800
$ws.Range("D48:M48").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 312900
$row[0,1] = 315100
$row[0,2] = 317500
$row[0,3] = 319900
$row[0,4] = 323700
$row[0,5] = 328700
$row[0,6] = 333600
$row[0,7] = 339800
$row[0,8] = 344900
$row[0,9] = 331900
$ws.Range("D49:M49").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D50:M50").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D51:M51").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 219700
$row[0,1] = 238200
$row[0,2] = 246700
$row[0,3] = 253400
$row[0,4] = 225500
$row[0,5] = 330700
$row[0,6] = 330300
$row[0,7] = 318900
$row[0,8] = 319000
$row[0,9] = 401300
$ws.Range("D52:M52").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D53:M53").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 4009500
$row[0,1] = 3866100
$row[0,2] = 3725100
$row[0,3] = 3683900
$row[0,4] = 3620100
$row[0,5] = 3927700
$row[0,6] = 3564100
$row[0,7] = 3516500
$row[0,8] = 3526200
$row[0,9] = 3242600
$ws.Range("D54:M54").Value = $row
$ws.Range("D55:M55").ClearContents() | Out-Null
$ws.Range("D56:M56").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 134800
$row[0,1] = 137000
$row[0,2] = 29800
$row[0,3] = 21000
$row[0,4] = 21600
$row[0,5] = 21400
$row[0,6] = 33500
$row[0,7] = 28300
$row[0,8] = 30700
$row[0,9] = 26200
$ws.Range("D57:M57").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 33900
$row[0,1] = 34900
$row[0,2] = 34900
$row[0,3] = 34800
$row[0,4] = 30100
$row[0,5] = 278000
$row[0,6] = 71300
$row[0,7] = 45400
$row[0,8] = 46100
$row[0,9] = 47200
$ws.Range("D58:M58").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 165000
$row[0,1] = 167400
$row[0,2] = 250000
$row[0,3] = 207600
$row[0,4] = 193500
$row[0,5] = 159900
$row[0,6] = 163600
$row[0,7] = 185000
$row[0,8] = 170000
$row[0,9] = 140900
$ws.Range("D59:M59").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 333700
$row[0,1] = 339300
$row[0,2] = 314600
$row[0,3] = 263400
$row[0,4] = 245300
$row[0,5] = 459300
$row[0,6] = 268500
$row[0,7] = 258700
$row[0,8] = 246800
$row[0,9] = 214200
$ws.Range("D60:M60").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 2821500
$row[0,1] = 2674900
$row[0,2] = 2575800
$row[0,3] = 2581700
$row[0,4] = 2587700
$row[0,5] = 2587200
$row[0,6] = 2427700
$row[0,7] = 2324200
$row[0,8] = 2376200
$row[0,9] = 2389400
$ws.Range("D61:M61").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 37400
$row[0,1] = 37900
$row[0,2] = 37500
$row[0,3] = 134700
$row[0,4] = 155400
$row[0,5] = 293900
$row[0,6] = 286800
$row[0,7] = 270000
$row[0,8] = 269800
$row[0,9] = 54500
$ws.Range("D62:M62").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D63:M63").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D64:M64").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D65:M65").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 3489900
$row[0,1] = 3349400
$row[0,2] = 3218300
$row[0,3] = 3257300
$row[0,4] = 3241400
$row[0,5] = 3585700
$row[0,6] = 3228200
$row[0,7] = 3148900
$row[0,8] = 3176400
$row[0,9] = 3052700
$ws.Range("D66:M66").Value = $row
$ws.Range("D67:M67").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D68:M68").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D69:M69").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D70:M70").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D71:M71").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 155900
$row[0,1] = 153900
$row[0,2] = 146200
$row[0,3] = 70400
$row[0,4] = 26100
$row[0,5] = 3100
$row[0,6] = -1900
$row[0,7] = 30800
$row[0,8] = 17600
$row[0,9] = 3800
$ws.Range("D72:M72").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D73:M73").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D74:M74").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D75:M75").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 519600
$row[0,1] = 516700
$row[0,2] = 506800
$row[0,3] = 426600
$row[0,4] = 378700
$row[0,5] = 342000
$row[0,6] = 335800
$row[0,7] = 367600
$row[0,8] = 349700
$row[0,9] = 190000
$ws.Range("D76:M76").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D77:M77").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 43465
$row[0,1] = 43373
$row[0,2] = 43281
$row[0,3] = 43190
$row[0,4] = 43100
$row[0,5] = 43008
$row[0,6] = 42916
$row[0,7] = 42825
$row[0,8] = 42735
$row[0,9] = 42643
$ws.Range("D80:M80").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 8900
$row[0,1] = 14700
$row[0,2] = 82700
$row[0,3] = 51200
$row[0,4] = 29500
$row[0,5] = 11800
$row[0,6] = -25700
$row[0,7] = 19900
$row[0,8] = 20400
$row[0,9] = 8300
$ws.Range("D81:M81").Value = $row
$ws.Range("D82:M82").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 46900
$row[0,1] = 44200
$row[0,2] = 46000
$row[0,3] = 43200
$row[0,4] = 43500
$row[0,5] = 42700
$row[0,6] = 46800
$row[0,7] = 45300
$row[0,8] = 42600
$row[0,9] = 36200
$ws.Range("D83:M83").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D84:M84").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D85:M85").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D86:M86").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D87:M87").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D88:M88").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 85000
$row[0,1] = 70000
$row[0,2] = 85000
$row[0,3] = 106000
$row[0,4] = 108400
$row[0,5] = 78600
$row[0,6] = 8100
$row[0,7] = 94900
$row[0,8] = 107200
$row[0,9] = 73100
$ws.Range("D89:M89").Value = $row
$ws.Range("D90:M90").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 371500
$row[0,1] = -149900
$row[0,2] = -120000
$row[0,3] = -137700
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 17900
$row[0,7] = -82700
$row[0,8] = -42900
$row[0,9] = -31900
$ws.Range("D91:M91").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D92:M92").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D93:M93").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = -202300
$row[0,1] = -152800
$row[0,2] = -114400
$row[0,3] = -137100
$row[0,4] = -82300
$row[0,5] = -66700
$row[0,6] = -90700
$row[0,7] = -42000
$row[0,8] = -34600
$row[0,9] = -326700
$ws.Range("D94:M94").Value = $row
$ws.Range("D95:M95").ClearContents() | Out-Null
$row = New-Object 'object[,]' 1,10
$row[0,0] = 20800
$row[0,1] = -6900
$row[0,2] = -6900
$row[0,3] = -6900
$row[0,4] = 0
$row[0,5] = 13300
$row[0,6] = -6800
$row[0,7] = -6600
$row[0,8] = -6500
$row[0,9] = 0
$ws.Range("D96:M96").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D97:M97").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D98:M98").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D99:M99").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 121300
$row[0,1] = 85200
$row[0,2] = -41200
$row[0,3] = -21100
$row[0,4] = -278200
$row[0,5] = 346800
$row[0,6] = 88800
$row[0,7] = -67100
$row[0,8] = -35100
$row[0,9] = 98600
$ws.Range("D100:M100").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D101:M101").Value = $row
$row = New-Object 'object[,]' 1,10
$row[0,0] = 4100
$row[0,1] = 2300
$row[0,2] = -70700
$row[0,3] = -52200
$row[0,4] = -252100
$row[0,5] = 358700
$row[0,6] = 6100
$row[0,7] = -14100
$row[0,8] = 37400
$row[0,9] = -155100
$ws.Range("D102:M102").Value = $row

Write-Output "done"